$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    14 = @(436, 393)
    15 = @(581, 398)
    16 = @(697, 429)
    17 = @(711, 383)
    18 = @(776, 383)
    19 = @(787, 365)
    20 = @(846, 368)
    21 = @(849, 482)
    22 = @(415, 476)
    23 = @(270, 477)
}

foreach ($row in $data.Keys | Sort-Object) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}

$ws.Range("B14").Select()
